$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '26.233.43'
$ws.Range("E2").Value = '  -0.91%  '

# Row 3
$ws.Range("D3").Value = '1.663.13'
$ws.Range("E3").Value = '  -0.82%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.005'
$ws.Range("E4").Value = '  +0.44%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '218.34'
$ws.Range("E5").Value = '  +0.41%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5219'
$ws.Range("E6").Value = '  -1.85%  '

# Row 7
$ws.Range("E7").Value = '  +0.43%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2648'
$ws.Range("E8").Value = '  -1.67%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06287'
$ws.Range("E9").Value = '  -2.03%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.85'
$ws.Range("E10").Value = '  -4.52%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07726'
$ws.Range("E11").Value = '  -1.31%  '

# Row 12
$ws.Range("D12").Value = '1.660.66'
$ws.Range("E12").Value = '  -1.03%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.432'
$ws.Range("E13").Value = '  -1.73%  '

# Row 14
$ws.Range("D14").Value = '1.889.31'
$ws.Range("E14").Value = '  -0.89%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5446'
$ws.Range("E15").Value = '  -2.41%  '

# Row 16
$ws.Range("D16").Value = '0.0₅8155'
$ws.Range("E16").Value = '  -2.23%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '64.60'
$ws.Range("E17").Value = '  -1.86%  '

# Row 18
$ws.Range("D18").Value = '26.251.13'

# Row 19
$ws.Range("E19").Value = '  +0.45%  '

# Row 20
$ws.Range("E20").Value = '  -2.59%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '193.33'
$ws.Range("E21").Value = '  -0.48%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.08'
$ws.Range("E22").Value = '  -2.14%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.051'
$ws.Range("E23").Value = '  -4.77%  '

# Row 24
$ws.Range("E24").Value = '  +0.58%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '139.98'
$ws.Range("E25").Value = '  -1.69%  '

# Row 26
$ws.Range("E26").Value = '  -4.58%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.167'
$ws.Range("E27").Value = '  -3.21%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '16.12'
$ws.Range("E28").Value = '  -1.12%  '

# Row 29
$ws.Range("E29").Value = '  -2.17%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.06093'
$ws.Range("E30").Value = '  -3.58%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.279'
$ws.Range("E31").Value = '  +0.44%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.583'
$ws.Range("E32").Value = '  -1.50%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.262'
$ws.Range("E33").Value = '  -5.66%  '

# Row 34
$ws.Range("E34").Value = '  -3.59%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9665'
$ws.Range("E35").Value = '  -4.40%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.429'
$ws.Range("E36").Value = '  +0.36%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.784'
$ws.Range("E37").Value = '  -0.15%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5675'
$ws.Range("E38").Value = '  -8.66%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01592'

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.983'
$ws.Range("E40").Value = '  -3.14%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8556'
$ws.Range("E41").Value = '  -1.24%  '

# Row 42
$ws.Range("E42").Value = '  +0.51%  '

# Row 43
$ws.Range("D43").Value = '1.012.98'
$ws.Range("E43").Value = '  -7.51%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '100.36'
$ws.Range("E44").Value = '  -0.17%  '

# Row 45
$ws.Range("D45").Value = '1.804.72'
$ws.Range("E45").Value = '  -0.96%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '57.10'
$ws.Range("E46").Value = '  -1.74%  '

# Row 47
$ws.Range("E47").Value = '  +3.39%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.009'
$ws.Range("E48").Value = '  +0.68%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.972'
$ws.Range("E49").Value = '  -2.80%  '

# Row 50
$ws.Range("B50").Value = 'RenderToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.480'
$ws.Range("E50").Value = '  -0.96%  '

# Row 51
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05182'
$ws.Range("E51").Value = '  -0.51%  '
